$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (was the old row 2 "heuristic" scenario, now gets the values
# that used to be on row 6 -- delays/tests switched on)
$ws.Range("D2").Value = 1115970.9
$ws.Range("F2").Value = 30000
$ws.Range("H2").Value = 15291647603.6941
$ws.Range("I2").Value = 102.9631116574518
$ws.Range("J2").Value = 15171371109.45141

# Update row 3 values (was the old row 3 "linearization_heuristic_Prop_Bouncing" scenario)
$ws.Range("D3").Value = 1115970.9
$ws.Range("F3").Value = 30000
$ws.Range("H3").Value = 15291647603.6941
$ws.Range("I3").Value = 102.9631116574518
$ws.Range("J3").Value = 15171371109.45141

# Remove rows 4 through 9 which held the now-deleted 30-day linearization heuristic benchmarks
$ws.Range("A4:J9").EntireRow.Delete()
